$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values (e.g. "1.001")
# are not auto-converted to numbers by Excel, matching the original inlineStr format.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "31.254.39"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "1.945.72"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "242.84"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.4829"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").Value = "0.2920"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").Value = "0.06820"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "20.18"
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("D11").Value = "104.94"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07836"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.956.15"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").Value = "5.326"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "0.6923"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "298.44"
$ws.Range("E16").Value = "  +8.57%  "
$ws.Range("D17").Value = "31.328.36"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").Value = "2.213.52"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").Value = "13.07"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "0.000007632"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "5.598"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("B23").Value = "BitDAO"
$ws.Range("C23").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D23").Value = "0.4733"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "6.476"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.601"
$ws.Range("E26").Value = "  -2.48%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "168.85"
$ws.Range("E27").Value = "  +2.36%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "19.89"
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "2.147"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "1.396"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.1018"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.642"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "1.540"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "4.372"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.04853"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.7457"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "1.138"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "2.739"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01969"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.653"
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.654"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "77.50"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "2.040"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8767"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.4386"
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "106.52"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "0.9993"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.026.94"
$ws.Range("E48").Value = "  +4.77%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "7.618"
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.219"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.1218"
$ws.Range("E51").Value = "  -1.51%  "
